# Update "想去人数" (people interested) counts that changed between the
# previous data snapshot and the newly generated one (gh-pages output
# regenerated at commit 456a3b4).
#
# Sheet "展览" (exhibitions only) and sheet "全部类型" (all types) both list
# the same events, so both copies of column F need to be bumped.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" ---
$wsExhibit.Range("F4").Value  = 616
$wsExhibit.Range("F5").Value  = 144
$wsExhibit.Range("F6").Value  = 9284
$wsExhibit.Range("F9").Value  = 1188
$wsExhibit.Range("F10").Value = 1087
$wsExhibit.Range("F11").Value = 141
$wsExhibit.Range("F14").Value = 256
$wsExhibit.Range("F16").Value = 83
$wsExhibit.Range("F18").Value = 1224

# --- Sheet "全部类型" ---
$wsAll.Range("F6").Value  = 616
$wsAll.Range("F7").Value  = 144
$wsAll.Range("F8").Value  = 9284
$wsAll.Range("F11").Value = 1188
$wsAll.Range("F12").Value = 1087
$wsAll.Range("F13").Value = 141
$wsAll.Range("F16").Value = 256
$wsAll.Range("F18").Value = 83
$wsAll.Range("F20").Value = 1225
